# Mise a jour taches et avancement - Agile.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Clear out the old "RAF per sprint" helper table (old B11:C16) and the
#     old G6/G7 trailing values / old formula row (row 8) so we can rebuild
#     the sheet cleanly with the new layout.
$ws.Range("B8:G16").ClearContents()
$ws.Range("M5:N16").ClearContents()

# --- Rows 8-10 are brand new task rows: give them the same bordered style
#     ("s=2") already used by rows 5-7 by copying the formatting over.
$ws.Range("B7:G7").Copy()
$ws.Range("B8:G10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Main task table (rows 5-10) ---------------------------------------
$ws.Range("B5").Value = "#1"
$ws.Range("C5").Value = "Montée en compétence sur Jmonkey & installation"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 100

$ws.Range("B6").Value = "#2"
$ws.Range("C6").Value = "Création des tâches et graphes agiles"
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 3
$ws.Range("G6").ClearContents()

$ws.Range("B7").Value = "#3"
$ws.Range("C7").Value = "Ouvrir une fenêtre avec un monde  16x16 blocs"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").ClearContents()

$ws.Range("B8").Value = "#4"
$ws.Range("C8").Value = "Vue à la première personne + imposer gravité"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 2
$ws.Range("G8").ClearContents()

$ws.Range("B9").Value = "#5"
$ws.Range("C9").Value = "Intégration système de collisions (solide)"
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 2
$ws.Range("G9").ClearContents()

$ws.Range("B10").Value = "#6"
$ws.Range("C10").Value = "Poser des blocs et enlever des blocs"
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 2
$ws.Range("G10").ClearContents()

# --- Totals row -----------------------------------------------------------
$ws.Range("E11").Formula = "=SUM(E5:E10)"
$ws.Range("F11").Formula = "=SUM(F5:F10)"
$ws.Range("G11").Value = 90

# --- New "Sprint / RAF" helper table used by the burn-down chart ----------
$ws.Range("M6").Value = "Sprint"
$ws.Range("N6").Value = "RAF"

$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 100

$ws.Range("M8").Value = 2
$ws.Range("N8").Value = 80

$ws.Range("M9").Value = 3
$ws.Range("N9").Value = 60

$ws.Range("M10").Value = 4
$ws.Range("N10").Value = 30

$ws.Range("M11").Value = 5
$ws.Range("N11").Value = 0

# --- Column C got a bit wider to fit the new task descriptions -----------
$ws.Columns.Item(3).ColumnWidth = 46.3

# --- Update chart series source to point at the new helper table ---------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(,Feuil1!`$M`$7:`$M`$11,Feuil1!`$N`$7:`$N`$11,1)"

# --- Move the chart further to the right / up on the sheet ---------------
# Anchor "from": 0-based xdr:col 15 / xdr:row 4 (i.e. 1-based column 16 / row 5)
$fromCol = 15
$fromColOff = 210911
$fromRow = 4
$fromRowOff = 34698

# Anchor "to": 0-based xdr:col 18 / xdr:row 18 (i.e. 1-based column 19 / row 19)
$toCol = 18
$toColOff = 639536
$toRow = 18
$toRowOff = 110898

$sumColWidthFrom = 0
for ($c = 1; $c -le $fromCol; $c++) {
    $sumColWidthFrom += $ws.Cells.Item(1, $c).Width
}
$sumRowHeightFrom = 0
for ($r = 1; $r -le $fromRow; $r++) {
    $sumRowHeightFrom += $ws.Cells.Item($r, 1).Height
}

$sumColWidthTo = 0
for ($c = 1; $c -le $toCol; $c++) {
    $sumColWidthTo += $ws.Cells.Item(1, $c).Width
}
$sumRowHeightTo = 0
for ($r = 1; $r -le $toRow; $r++) {
    $sumRowHeightTo += $ws.Cells.Item($r, 1).Height
}

$newLeft = $sumColWidthFrom + ($fromColOff / 12700)
$newTop = $sumRowHeightFrom + ($fromRowOff / 12700)
$newRight = $sumColWidthTo + ($toColOff / 12700)
$newBottom = $sumRowHeightTo + ($toRowOff / 12700)

$co.Left = $newLeft
$co.Top = $newTop
$co.Width = $newRight - $newLeft
$co.Height = $newBottom - $newTop

# --- View housekeeping: zoom back to 100%, reset scroll position, select E9
$win = $excel.ActiveWindow
$win.Zoom = 100
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("E9").Select()
